$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the IP Address in A2 to the new value
$ws.Range("A2").Value = "192.168.146.135"

# Move the active selection to G2
$ws.Range("G2").Select()
